# Update average_county_temperature (column AD) values for specific facility
# rows to reflect updated NOAA temperature data, per the commit:
#   "Updated temperature with NOAA data"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Facilities H=1000606 (rows 2-5), H=1000615 & H=1001985 (rows 50-89),
# H=1002285 (rows 119-130)
$ws.Range("AD2:AD5").Value = 1.925925925925943
$ws.Range("AD50:AD89").Value = 1.925925925925943
$ws.Range("AD119:AD130").Value = 1.925925925925943

# Facility H=1002283 (rows 90-106)
$ws.Range("AD90:AD106").Value = -1.226851851851833

# Facility H=1004369 & H=1005361 (rows 175-206)
$ws.Range("AD175:AD206").Value = 13.17361111111111
